$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new manager rows below the existing data (rows 2-3 already present)
$ws.Range("A4").Value = "James"
$ws.Range("B4").Value = "T1212121C"
$ws.Range("C4").Value = 32
$ws.Range("D4").Value = "Married"
$ws.Range("E4").Value = "password"

$ws.Range("A5").Value = "Frank"
$ws.Range("B5").Value = "S2323232H"
$ws.Range("C5").Value = 30
$ws.Range("D5").Value = "Single"
$ws.Range("E5").Value = "password"

$ws.Range("A6").Value = "Kelly"
$ws.Range("B6").Value = "T9912834K"
$ws.Range("C6").Value = 44
$ws.Range("D6").Value = "Married"
$ws.Range("E6").Value = "password"

# Update the selected cell to match the final saved state
$ws.Range("G6").Select()
